$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "max limit" descriptor text values
$ws.Range("Q10").Value = "Pro32xD,Pro16xD,Pro815D,Pro885D"
$ws.Range("R10").Value = "MX 252,P485D,Pro215D,Pro32xBB,Pro16xBB"
$ws.Range("Q11").Value = "Pro32xBB,Pro16xBB"
$ws.Range("R11").Value = "MX 252,P485D,Pro215D,Pro815D,Pro885D"

# Row heights grow to fit the updated wrapped text
$ws.Rows.Item(10).RowHeight = 86.4
$ws.Rows.Item(11).RowHeight = 86.4

# Update selection to reflect where the editor left off
$ws.Range("I11").Select()
